$d = $word.ActiveDocument

# Paragraph 13 (1-based) is the first of the two consecutive empty
# paragraphs right after the "AOS.js (Animate On Scroll)" bullet in the
# "Animacje" section. It receives a plain run with no direct run
# formatting.
$p1 = $d.Paragraphs.Item(13)
$p1.Range.Text = "Aby użyć ikon Font Awesome, zainstaluj bibliotekę:"

# Paragraph 14 (1-based) is the next (also originally empty) paragraph.
# It receives a run styled like the rest of the body text: Times New
# Roman, 12pt (half-point size 24), with the Polish east-Asian language
# tag carried over from the paragraph mark's rPr.
$p2 = $d.Paragraphs.Item(14)
$r2 = $p2.Range
$r2.Text = "npm install @fortawesome/fontawesome-free"
$r2.Font.Name = "Times New Roman"
$r2.Font.NameAscii = "Times New Roman"
$r2.Font.NameFarEast = "Times New Roman"
$r2.Font.NameOther = "Times New Roman"
$r2.Font.NameBi = "Times New Roman"
$r2.Font.Size = 12
$r2.Font.SizeBi = 12
$r2.LanguageIDFarEast = "pl-PL"
